$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typos "then" -> "than" in the cgt/clt description rows (54, 55)
$ws.Range("D54").Value = "Push true to the stack if value 1 is greater than value 2, else push false to the stack"
$ws.Range("D55").Value = "Push true to the stack if value 1 is less than value 2, else push false to the stack"

# Row 56 used to hold the "syscall" opcode (0x36); it is being replaced by the
# new "cgteq" opcode, and syscall moves down to 0x38 (row 58)
$ws.Range("B56").Value = "cgteq"
$ws.Range("C56").ClearContents()
$ws.Range("D56").Value = "Push true to the stack if value 1 is greater than or equal to value 2, else push false to the stack"

# New row 57: clteq opcode (0x37)
$ws.Range("A57").Value = "0x37"
$ws.Range("B57").Value = "clteq"
$ws.Range("D57").Value = "Push true to the stack if value 1 is less than or equal to value 2, else push false to the stack"

# New row 58: syscall opcode, now 0x38
$ws.Range("A58").Value = "0x38"
$ws.Range("B58").Value = "syscall"
$ws.Range("C58").Value = "<uint16 (code)>"
$ws.Range("D58").Value = "Run a system function specified by (code), uses other values in the stack for arguments"

# Leave the final selection where the author left it in the commit
$ws.Range("D52").Select()
